$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 4) of job-search tracker data, mirroring the existing
# layout (row 3) but for a new application entry.
$ws.Range("B4").Value = "Pretend Corp."
$ws.Range("C4").Value = "Software Developer"
$ws.Range("D4").Value = "not yet"
$ws.Range("F4").Value = "Blaze Smith"
$ws.Range("G4").Value = "alexrogers823@yahoo.com"
$ws.Range("J4").Value = "Email Sent"
$ws.Range("K4").Value = "www.notimportant.com"
$ws.Range("L4").Value = "Open"
$ws.Range("M4").Value = "Python"
